# Generate Report for Handback
# Adds two newly handed-back files (b2891065-... and c9e26659-...) as new
# rows (6 & 7) across the "Overview", "zh-cn" and "de-de" worksheets,
# mirroring the existing rows' layout, hyperlinks and formatting.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"
$includeText = "Include"

# ---------------------------------------------------------------------
# New file identities
# ---------------------------------------------------------------------
$file1 = "b2891065-cfe7-4731-94c2-4c754d4bef1c"
$file1Hash = "e243ccac099635fd236a26a83151e14520daf803"
$file2 = "c9e26659-4cc2-4290-ae4b-5198713a52b8"
$file2Hash = "27a04f146edc21599d6f9ab6bfa524c441f7c242"

$file1Md = $file1 + ".md"
$file2Md = $file2 + ".md"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Cells.Item(6, 2).Value = $statusText
$wsOverview.Cells.Item(6, 3).Value = $statusText
$wsOverview.Cells.Item(7, 2).Value = $statusText
$wsOverview.Cells.Item(7, 3).Value = $statusText

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("A6"),
    "https://github.com/OpenLocalizationTest/oltest/blob/" + $file1Hash + "/e2e/" + $file1Md,
    [Type]::Missing,
    [Type]::Missing,
    $file1Md) | Out-Null

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("A7"),
    "https://github.com/OpenLocalizationTest/oltest/blob/" + $file2Hash + "/e2e/" + $file2Md,
    [Type]::Missing,
    [Type]::Missing,
    $file2Md) | Out-Null

# ---------------------------------------------------------------------
# Helper data driving the "zh-cn" / "de-de" detail sheets
# ---------------------------------------------------------------------
$langs = @(
    @{
        SheetName   = "zh-cn"
        Lang        = "zh-cn"
        HandoffOrg  = "OpenLocalizationTestOrg/oltest.zh-cn"
        Datetime1Off = "2016-02-22 18:04:37"
        Datetime1Back = "2016-02-22 18:05:22"
        Datetime2Off = "2016-02-22 18:04:37"
        Datetime2Back = "2016-02-22 18:05:22"
    },
    @{
        SheetName   = "de-de"
        Lang        = "de-de"
        HandoffOrg  = "OpenLocalizationTestOrg/oltest.de-de"
        Datetime1Off = "2016-02-22 18:04:49"
        Datetime1Back = "2016-02-22 18:05:42"
        Datetime2Off = "2016-02-22 18:04:49"
        Datetime2Back = "2016-02-22 18:05:42"
    }
)

foreach ($entry in $langs) {
    $ws = $wb.Worksheets.Item($entry.SheetName)
    $lang = $entry.Lang

    $xlf1 = $file1 + "." + $file1Hash + "." + $lang + ".xlf"
    $xlf2 = $file2 + "." + $file2Hash + "." + $lang + ".xlf"

    # ---- Row 6 : b2891065-... ------------------------------------------------
    $ws.Cells.Item(6, 2).Value = $statusText
    $ws.Cells.Item(6, 4).Value = $entry.Datetime1Off
    $ws.Cells.Item(6, 7).Value = $entry.Datetime1Back
    $ws.Cells.Item(6, 8).Value = $includeText

    $ws.Hyperlinks.Add(
        $ws.Range("A6"),
        "https://github.com/OpenLocalizationTest/oltest/blob/" + $file1Hash + "/e2e/" + $file1Md,
        [Type]::Missing,
        [Type]::Missing,
        $file1Md) | Out-Null

    $ws.Hyperlinks.Add(
        $ws.Range("C6"),
        "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/" + $file1Hash + "/ol-handoff/" + $entry.HandoffOrg + "/xinjiang/ht/" + $xlf1,
        [Type]::Missing,
        [Type]::Missing,
        $xlf1) | Out-Null

    $ws.Hyperlinks.Add(
        $ws.Range("E6"),
        "https://github.com/" + $entry.HandoffOrg + "/blob/" + $file1Hash + "/e2e/" + $file1Md,
        [Type]::Missing,
        [Type]::Missing,
        $file1Md) | Out-Null

    $ws.Hyperlinks.Add(
        $ws.Range("F6"),
        "https://github.com/OpenLocalizationTestOrg/olhandback/blob/" + $file1Hash + "/ol-handback/" + $entry.HandoffOrg + "/xinjiang/ht/" + $xlf1,
        [Type]::Missing,
        [Type]::Missing,
        $xlf1) | Out-Null

    # ---- Row 7 : c9e26659-... ------------------------------------------------
    $ws.Cells.Item(7, 2).Value = $statusText
    $ws.Cells.Item(7, 4).Value = $entry.Datetime2Off
    $ws.Cells.Item(7, 7).Value = $entry.Datetime2Back
    $ws.Cells.Item(7, 8).Value = $includeText

    $ws.Hyperlinks.Add(
        $ws.Range("A7"),
        "https://github.com/OpenLocalizationTest/oltest/blob/" + $file2Hash + "/e2e/" + $file2Md,
        [Type]::Missing,
        [Type]::Missing,
        $file2Md) | Out-Null

    $ws.Hyperlinks.Add(
        $ws.Range("C7"),
        "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/" + $file2Hash + "/ol-handoff/" + $entry.HandoffOrg + "/xinjiang/ht/" + $xlf2,
        [Type]::Missing,
        [Type]::Missing,
        $xlf2) | Out-Null

    $ws.Hyperlinks.Add(
        $ws.Range("E7"),
        "https://github.com/" + $entry.HandoffOrg + "/blob/" + $file2Hash + "/e2e/" + $file2Md,
        [Type]::Missing,
        [Type]::Missing,
        $file2Md) | Out-Null

    $ws.Hyperlinks.Add(
        $ws.Range("F7"),
        "https://github.com/OpenLocalizationTestOrg/olhandback/blob/" + $file2Hash + "/ol-handback/" + $entry.HandoffOrg + "/xinjiang/ht/" + $xlf2,
        [Type]::Missing,
        [Type]::Missing,
        $xlf2) | Out-Null

    # Match the date-time number format used by the existing rows.
    $ws.Range("D6:D7").NumberFormat = "yyyy-mm-dd HH:mm:ss"
    $ws.Range("G6:G7").NumberFormat = "yyyy-mm-dd HH:mm:ss"
}

Write-Host "Handback rows added."
